$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1350
$ws.Range("F4").Value = 1156
$ws.Range("F5").Value = 1054
$ws.Range("F6").Value = 1842
$ws.Range("F7").Value = 585
$ws.Range("F8").Value = 1226
$ws.Range("F12").Value = 315
$ws.Range("F13").Value = 91
$ws.Range("F15").Value = 736
$ws.Range("F16").Value = 199
$ws.Range("F17").Value = 112
$ws.Range("F22").Value = 685
$ws.Range("F23").Value = 54
$ws.Range("F25").Value = 172

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 18
$ws.Range("F11").Value = 124

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 318

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 318
$ws.Range("F4").Value = 1350
$ws.Range("F5").Value = 1156
$ws.Range("F6").Value = 1054
$ws.Range("F7").Value = 1842
$ws.Range("F8").Value = 585
$ws.Range("F9").Value = 1226
$ws.Range("F14").Value = 315
$ws.Range("F15").Value = 91
$ws.Range("F17").Value = 736
$ws.Range("F18").Value = 199
$ws.Range("F19").Value = 112
$ws.Range("F24").Value = 18
$ws.Range("F30").Value = 685
$ws.Range("F31").Value = 54
$ws.Range("F33").Value = 172
$ws.Range("F43").Value = 124
$ws.Range("F44").Value = 124
